$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Type" column (J, values "Pool", data validation "Pool,All") is being
# removed entirely. Column K ("Rule For", Accounting/Reporting validation)
# shifts left to become the new column J.

# Capture the text of column K's header comment ("Rule For") before the
# structural edit so it can be carried over to (the new) J1 once column J
# is removed. Replacing J1's existing comment in place (rather than
# deleting then re-adding) keeps it on the original comment author and
# lets us restore the bold "Author:" lead-in run.
$k1Comment = $ws.Range("K1").Comment
$k1Text = $k1Comment.Text()

$newJ1Comment = $ws.Range("J1").AddComment($k1Text)
$newJ1Comment.Shape.TextFrame.Characters(1, 7).Font.Bold = $true

# Drop the now-duplicated comment sitting on K1.
$ws.Range("K1").Comment.Delete()

# Delete column J (the "Type" column). This shifts K -> J, updates the
# dimension, row spans, data validation range, and carries the comment we
# just re-homed on J1 along with the shift.
$ws.Columns("J").Delete()

# Match the author's final selection (whole column J, as if the new
# "Rule For" column / its validation dropdown was selected).
$ws.Range("J1:J1048576").Select() | Out-Null
